$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: communication -> category_shopping (merged with shopping benefits + extras)
$ws.Range("A2").Value = "category_shopping"
$ws.Range("B2").Value = "백화점,마트/편의점,온라인쇼핑,소셜커머스,해외직구,홈쇼핑,SPA브랜드,아울렛,대형마트,SSM,전통시장,면세점,모든가맹점,국내외가맹점,화장품,드럭스토어"

# Row 3: shopping -> category_culture (culture benefits + golf)
$ws.Range("A3").Value = "category_culture"
$ws.Range("B3").Value = "영화,공연/전시,문화센터,도서,음원사이트,영화/문화,디지털구독,골프"

# Row 4: culture -> category_transportation
$ws.Range("A4").Value = "category_transportation"
$ws.Range("B4").Value = "고속버스,항공권,기차,대중교통,렌터카,택시,교통"

# Row 5: transportation -> category_car
$ws.Range("A5").Value = "category_car"
$ws.Range("B5").Value = "주유,자동차,정비,주유소,충전소,하이패스,자동차/하이패스"

# Row 6: food -> category_food (benefit text unchanged)
$ws.Range("A6").Value = "category_food"

# Row 7: education -> category_education (benefit text unchanged)
$ws.Range("A7").Value = "category_education"

# Row 8: utilities -> category_housing_communication
$ws.Range("A8").Value = "category_housing_communication"
$ws.Range("B8").Value = "공과금,공과금/렌탈,통신,KT,SKT,LGU+"

# Row 9: aviation -> category_travel
$ws.Range("A9").Value = "category_travel"
$ws.Range("B9").Value = "진에어,대한항공,아시아나항공,제주항공,항공마일리지,온라인 여행사,여행/숙박,공항라운지,리조트,공항,공항라운지/PP,여행사,호텔,PP,라운지키,해외이용,저가항공,해외"

# Row 10: medical -> category_medical
$ws.Range("A10").Value = "category_medical"
$ws.Range("B10").Value = "약국,병원,병원/약국"

# Row 11: others -> category_financial_insurance
$ws.Range("A11").Value = "category_financial_insurance"
$ws.Range("B11").Value = "금융,은행사,보험,보험사"

# Row 12 (new): category_others
$ws.Range("A12").Value = "category_others"
$ws.Range("B12").Value = "애완동물,비즈니스,동물병원,지역,생활,렌탈,무실적"

$wb.Save()
